$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B13").Value = "Inital Solution 0 changes  totQ=800.00 totD=138.99 `$138.99"
$ws.Range("C13").Value = "Inital Solution 0 changes  totQ=800.00 totD=125.02 `$125.02"

$ws.Range("T13").Value = "Inital Solution 0 changes  totQ=800.00 totD=110.97 `$110.97"
$ws.Range("U13").Value = "Inital Solution 0 changes  totQ=800.00 totD=110.97 `$110.97"
$ws.Range("V13").Value = "Inital Solution 0 changes  totQ=800.00 totD=106.78 `$106.78"
$ws.Range("W13").Value = "Inital Solution 0 changes  totQ=800.00 totD=106.78 `$106.78"
$ws.Range("X13").Value = "Inital Solution 0 changes  totQ=800.00 totD=114.96 `$114.96"
$ws.Range("Y13").Value = "Inital Solution 0 changes  totQ=800.00 totD=114.96 `$114.96"
